$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update existing row 2 scores ----
$ws.Range("B2").Value = 75
$ws.Range("C2").Value = 70
$ws.Range("D2").Value = 70
$ws.Range("E2").Value = 60
$ws.Range("F2").Value = 55
$ws.Range("G2").Value = 79

# ---- Apply the same style as rows 2/3 to the new block of rows 4-31 ----
$ws.Range("A2:G3").Copy()
$ws.Range("A4:G31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Add new student rows (John, Alice) with their scores ----
$ws.Range("A4").Value = "John"
$ws.Range("B4").Value = 69
$ws.Range("C4").Value = 28
$ws.Range("D4").Value = 32
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = 26
$ws.Range("G4").Value = 87

$ws.Range("A5").Value = "Alice"
$ws.Range("B5").Value = 13
$ws.Range("C5").Value = 33
$ws.Range("D5").Value = 22
$ws.Range("E5").Value = 19
$ws.Range("F5").Value = 67
$ws.Range("G5").Value = 32

# ---- Rows 6-31 stay blank (formatted, but empty) ----

# ---- Update the selected/active range to mirror the report view ----
$ws.Range("A6:G31").Select()

Write-Output "edit applied"
